$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.613.60'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '2.668.03'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '597.61'
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").Value = '157.45'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.613'
$ws.Range("E8").Value = '  +4.18%  '
$ws.Range("D9").Value = '0.128'
$ws.Range("E9").Value = '  +2.28%  '
$ws.Range("D10").Value = '0.398'
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("D11").Value = '5.82'
$ws.Range("E11").Value = '  -3.08%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '29.01'
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").Value = '0.0000198'
$ws.Range("E14").Value = '  -4.47%  '
$ws.Range("D15").Value = '3.151.56'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '65.571.57'
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").Value = '2.677.50'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = '12.65'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").Value = '4.78'
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("D20").Value = '7.49'
$ws.Range("E20").Value = '  -3.00%  '
$ws.Range("D21").Value = '351.19'
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").Value = '0.0000113'
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").Value = '9.71'
$ws.Range("E25").Value = '  -2.49%  '
$ws.Range("D26").Value = '1.67'
$ws.Range("E26").Value = '  +2.66%  '
$ws.Range("E27").Value = '  -3.43%  '
$ws.Range("D28").Value = '0.165'
$ws.Range("E28").Value = '  -3.62%  '
$ws.Range("D29").Value = '7.98'
$ws.Range("E29").Value = '  -3.29%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.13'
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '530.28'
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").Value = '1.78'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = '6.46'
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("D35").Value = '5.44'
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").Value = '0.423'
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("D37").Value = '20.55'
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").Value = '156.82'
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("D40").Value = '1.93'
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '162.60'
$ws.Range("E42").Value = '  -3.06%  '
$ws.Range("D43").Value = '4.09'
$ws.Range("E43").Value = '  -1.38%  '
$ws.Range("D44").Value = '2.33'
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("D45").Value = '0.0609'
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("D46").Value = '22.62'
$ws.Range("E46").Value = '  -4.27%  '
$ws.Range("D47").Value = '0.640'
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("E48").Value = '  -3.12%  '
$ws.Range("D49").Value = '0.0₆0252'
$ws.Range("E49").Value = '  +6.82%  '
$ws.Range("D50").Value = '0.0985'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").Value = '19.82'
$ws.Range("E51").Value = '  -4.48%  '
